$d = $word.ActiveDocument

# --- Paragraph 1 ("ABOUT ME"): switch font to Baskerville Old Face ---
$p1 = $d.Paragraphs(1)
$p1.Range.Font.Name = "Baskerville Old Face"

# The _GoBack bookmark currently sits in paragraph 1; it needs to move into
# paragraph 2 (between the two sentences of the rewritten text), so drop it
# here and re-add it later at the right spot.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Paragraph 2: rewrite the "about me" blurb and switch font to Comic Sans MS ---
$p2 = $d.Paragraphs(2)
$p2Start = $p2.Range.Start

$firstPart  = "I am passionately driven to provide quality services with dynamism in an environment that nurtures and optimizes skills, contributing effectively to organizational success. With a consistently positive attitude towards "
$secondPart = "work and a can-do spirit essential for any team, I have extensive experience in business administration. Now, I am eager to transition into the tech industry as I undergo training with Testify Ltd. My expectations are high as I have collaborated with exceptional individuals who have shared their knowledge with me, and I am confident in securing a software testing QA role soon."
$newText = $firstPart + $secondPart

$p2.Range.Text = $newText
$p2.Range.Font.Name = "Comic Sans MS"

# Re-insert _GoBack exactly at the boundary between the two parts of the text.
$splitOffset = $p2Start + $firstPart.Length
$bmRange = $d.Range($splitOffset, $splitOffset)
$d.Bookmarks.Add("_GoBack", $bmRange)
